$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'768"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "'1743687.79"
$ws.Range("D2").Style = "Normal"

$ws.Range("C4").Value = "'1023"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "'3635453.47"
$ws.Range("D4").Style = "Normal"

$ws.Range("C6").Value = "'664"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'2152907.78"
$ws.Range("D6").Style = "Normal"

$ws.Range("C9").Value = "'193"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'570715.60"
$ws.Range("D9").Style = "Normal"

$ws.Range("C10").Value = "'376"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'1414220.18"
$ws.Range("D10").Style = "Normal"

$ws.Range("C11").Value = "'174"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'640976.34"
$ws.Range("D11").Style = "Normal"

$ws.Range("C14").Value = "'223"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'602362.00"
$ws.Range("D14").Style = "Normal"

$ws.Range("C15").Value = "'7"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'22500.00"
$ws.Range("D15").Style = "Normal"

$ws.Range("C16").Value = "'498"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'1846574.13"
$ws.Range("D16").Style = "Normal"

$ws.Range("C17").Value = "'143"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'445375.81"
$ws.Range("D17").Style = "Normal"

$ws.Range("C20").Value = "'182"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'460659.00"
$ws.Range("D20").Style = "Normal"

$ws.Range("C21").Value = "'337"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'1195741.00"
$ws.Range("D21").Style = "Normal"

$ws.Range("C30").Value = "'565"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'2338016.89"
$ws.Range("D30").Style = "Normal"

$ws.Range("C32").Value = "'393"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'1347198.57"
$ws.Range("D32").Style = "Normal"

$ws.Range("C45").Value = "'394"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'1061299.43"
$ws.Range("D45").Style = "Normal"

$ws.Range("C47").Value = "'629"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'2512013.99"
$ws.Range("D47").Style = "Normal"

$ws.Range("C48").Value = "'430"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'1502861.40"
$ws.Range("D48").Style = "Normal"

$ws.Range("C61").Value = "'24"
$ws.Range("C61").Style = "Normal"
$ws.Range("D61").Value = "'84777.00"
$ws.Range("D61").Style = "Normal"
